$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 12 ("Enterprises (absolute #)" / "40000")
# and row 13 ("Enterprises density (per 1000 people)" / "38.7") so that
# the density row now appears first, followed by the absolute # row.

$a12 = $ws.Range("A12").Value()
$d12 = $ws.Range("D12").Value()
$a13 = $ws.Range("A13").Value()
$d13 = $ws.Range("D13").Value()

$ws.Range("A12").Value = $a13
# Prefix numeric-looking text with an apostrophe so Excel keeps storing it
# as text (matching the original shared-string cell type) instead of
# silently converting it to a number.
$ws.Range("D12").Value = "'" + $d13
$ws.Range("A13").Value = $a12
$ws.Range("D13").Value = "'" + $d12
